$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1965.3478
$ws.Cells.Item(40, 9).Value = 1323.1538
$ws.Cells.Item(40, 10).Value = 2800.2
$ws.Cells.Item(40, 11).Value = 1323.1538
$ws.Cells.Item(40, 12).Value = 2800.2
$ws.Cells.Item(40, 13).Value = -1148.1538
$ws.Cells.Item(40, 14).Value = -3150.2

$ws.Cells.Item(70, 8).Value = 1409.5238
$ws.Cells.Item(70, 9).Value = 1709.1818
$ws.Cells.Item(70, 10).Value = 1079.9
$ws.Cells.Item(70, 11).Value = 5127.5454
$ws.Cells.Item(70, 12).Value = 3239.7
$ws.Cells.Item(70, 13).Value = -4857.5454
$ws.Cells.Item(70, 14).Value = -3779.7

$ws.Cells.Item(73, 8).Value = 1409.5238
$ws.Cells.Item(73, 9).Value = 1709.1818
$ws.Cells.Item(73, 10).Value = 1079.9
$ws.Cells.Item(73, 11).Value = 5127.5454
$ws.Cells.Item(73, 12).Value = 3239.7
$ws.Cells.Item(73, 13).Value = -4191.5454
$ws.Cells.Item(73, 14).Value = -5111.700000000001

$ws.Cells.Item(74, 8).Value = 3247.0952
$ws.Cells.Item(74, 9).Value = 3141.6667
$ws.Cells.Item(74, 10).Value = 3387.6667
$ws.Cells.Item(74, 11).Value = 3141.6667
$ws.Cells.Item(74, 12).Value = 3387.6667
$ws.Cells.Item(74, 13).Value = -2205.6667
$ws.Cells.Item(74, 14).Value = -5259.6667

$ws.Cells.Item(76, 8).Value = 4676.8696
$ws.Cells.Item(76, 9).Value = 4616.4707
$ws.Cells.Item(76, 10).Value = 4848
$ws.Cells.Item(76, 11).Value = 4616.4707
$ws.Cells.Item(76, 12).Value = 4848
$ws.Cells.Item(76, 13).Value = -4301.4707
$ws.Cells.Item(76, 14).Value = -5478

$ws.Cells.Item(77, 8).Value = 3247.0952
$ws.Cells.Item(77, 9).Value = 3141.6667
$ws.Cells.Item(77, 10).Value = 3387.6667
$ws.Cells.Item(77, 11).Value = 15708.3335
$ws.Cells.Item(77, 12).Value = 16938.3335
$ws.Cells.Item(77, 13).Value = -11028.3335
$ws.Cells.Item(77, 14).Value = -26298.3335

$ws.Cells.Item(79, 8).Value = 4676.8696
$ws.Cells.Item(79, 9).Value = 4616.4707
$ws.Cells.Item(79, 10).Value = 4848
$ws.Cells.Item(79, 11).Value = 4616.4707
$ws.Cells.Item(79, 12).Value = 4848
$ws.Cells.Item(79, 13).Value = -3524.4707
$ws.Cells.Item(79, 14).Value = -7032

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7021.51
$ws.Cells.Item(32, 9).Value = 5095.9214
$ws.Cells.Item(32, 10).Value = 22601.273
$ws.Cells.Item(32, 11).Value = 5095.9214
$ws.Cells.Item(32, 12).Value = 22601.273
$ws.Cells.Item(32, 13).Value = -4808.9214
$ws.Cells.Item(32, 14).Value = -23175.273

$ws.Cells.Item(45, 8).Value = 131955.25
$ws.Cells.Item(45, 9).Value = 251500.75
$ws.Cells.Item(45, 10).Value = 12409.75
$ws.Cells.Item(45, 11).Value = 251500.75
$ws.Cells.Item(45, 12).Value = 12409.75
$ws.Cells.Item(45, 13).Value = -251123.75
$ws.Cells.Item(45, 14).Value = -13163.75

$ws.Cells.Item(63, 8).Value = 2397.4614
$ws.Cells.Item(63, 9).Value = 2150.8333
$ws.Cells.Item(63, 10).Value = 2608.8572
$ws.Cells.Item(63, 11).Value = 2150.8333
$ws.Cells.Item(63, 12).Value = 2608.8572
$ws.Cells.Item(63, 13).Value = -1464.8333
$ws.Cells.Item(63, 14).Value = -3980.8572

$ws.Cells.Item(66, 8).Value = 2397.4614
$ws.Cells.Item(66, 9).Value = 2150.8333
$ws.Cells.Item(66, 10).Value = 2608.8572
$ws.Cells.Item(66, 11).Value = 10754.1665
$ws.Cells.Item(66, 12).Value = 13044.286
$ws.Cells.Item(66, 13).Value = -7322.166499999999
$ws.Cells.Item(66, 14).Value = -19908.286

$ws.Cells.Item(74, 8).Value = 1193.0322
$ws.Cells.Item(74, 9).Value = 1258.5
$ws.Cells.Item(74, 10).Value = 968.5714
$ws.Cells.Item(74, 11).Value = 1258.5
$ws.Cells.Item(74, 12).Value = 968.5714
$ws.Cells.Item(74, 13).Value = -384.5
$ws.Cells.Item(74, 14).Value = -2716.5714

$ws.Cells.Item(77, 8).Value = 1193.0322
$ws.Cells.Item(77, 9).Value = 1258.5
$ws.Cells.Item(77, 10).Value = 968.5714
$ws.Cells.Item(77, 11).Value = 6292.5
$ws.Cells.Item(77, 12).Value = 4842.857
$ws.Cells.Item(77, 13).Value = -1924.5
$ws.Cells.Item(77, 14).Value = -13578.857

$ws.Cells.Item(106, 8).Value = 44498
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 44498
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 44498
$ws.Cells.Item(106, 14).Value = -47022

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 183751
$ws.Cells.Item(105, 9).Value = 144981.42
$ws.Cells.Item(105, 10).Value = 251597.75
$ws.Cells.Item(105, 11).Value = 144981.42
$ws.Cells.Item(105, 12).Value = 251597.75
$ws.Cells.Item(105, 13).Value = -143234.42
$ws.Cells.Item(105, 14).Value = -255091.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2507.5334
$ws.Cells.Item(58, 9).Value = 2427.182
$ws.Cells.Item(58, 10).Value = 2728.5
$ws.Cells.Item(58, 11).Value = 2427.182
$ws.Cells.Item(58, 12).Value = 2728.5
$ws.Cells.Item(58, 13).Value = -2224.182
$ws.Cells.Item(58, 14).Value = -3134.5

$ws.Cells.Item(122, 8).Value = 993.3333
$ws.Cells.Item(122, 9).Value = 993.3333
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 2979.9999
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -529.9998999999998
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 2507.5334
$ws.Cells.Item(136, 9).Value = 2427.182
$ws.Cells.Item(136, 10).Value = 2728.5
$ws.Cells.Item(136, 11).Value = 7281.545999999999
$ws.Cells.Item(136, 12).Value = 8185.5
$ws.Cells.Item(136, 13).Value = -4731.545999999999
$ws.Cells.Item(136, 14).Value = -13285.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 932.92
$ws.Cells.Item(113, 9).Value = 1748.75
$ws.Cells.Item(113, 10).Value = 549
$ws.Cells.Item(113, 11).Value = 5246.25
$ws.Cells.Item(113, 12).Value = 1647
$ws.Cells.Item(113, 13).Value = -3076.25
$ws.Cells.Item(113, 14).Value = -5987

$ws.Cells.Item(122, 8).Value = 526.087
$ws.Cells.Item(122, 9).Value = 599
$ws.Cells.Item(122, 10).Value = 515.15
$ws.Cells.Item(122, 11).Value = 5391
$ws.Cells.Item(122, 12).Value = 4636.349999999999
$ws.Cells.Item(122, 13).Value = -2941
$ws.Cells.Item(122, 14).Value = -9536.349999999999

$ws.Cells.Item(140, 8).Value = 6629.1577
$ws.Cells.Item(140, 9).Value = 6629.1577
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 19887.4731
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -14707.4731
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 96648.91
$ws.Cells.Item(70, 9).Value = 138183.73
$ws.Cells.Item(70, 10).Value = 7645.7144
$ws.Cells.Item(70, 11).Value = 138183.73
$ws.Cells.Item(70, 12).Value = 7645.7144
$ws.Cells.Item(70, 13).Value = -137913.73
$ws.Cells.Item(70, 14).Value = -8185.7144

$ws.Cells.Item(73, 8).Value = 96648.91
$ws.Cells.Item(73, 9).Value = 138183.73
$ws.Cells.Item(73, 10).Value = 7645.7144
$ws.Cells.Item(73, 11).Value = 138183.73
$ws.Cells.Item(73, 12).Value = 7645.7144
$ws.Cells.Item(73, 13).Value = -137247.73
$ws.Cells.Item(73, 14).Value = -9517.714400000001

$ws.Cells.Item(80, 8).Value = 250258750
$ws.Cells.Item(80, 9).Value = 250258750
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 250258750
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -250257752
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(83, 8).Value = 250258750
$ws.Cells.Item(83, 9).Value = 250258750
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 1251293750
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -1251288758
$ws.Cells.Item(83, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 1882.4
$ws.Cells.Item(113, 9).Value = 1800
$ws.Cells.Item(113, 10).Value = 1891.5555
$ws.Cells.Item(113, 11).Value = 1800
$ws.Cells.Item(113, 12).Value = 1891.5555
$ws.Cells.Item(113, 13).Value = 370
$ws.Cells.Item(113, 14).Value = -6231.5555

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 23079.914
$ws.Cells.Item(40, 9).Value = 37130.57
$ws.Cells.Item(40, 10).Value = 2373.6843
$ws.Cells.Item(40, 11).Value = 37130.57
$ws.Cells.Item(40, 12).Value = 2373.6843
$ws.Cells.Item(40, 13).Value = -36994.57
$ws.Cells.Item(40, 14).Value = -2645.6843

$ws.Cells.Item(136, 8).Value = 1553.3414
$ws.Cells.Item(136, 9).Value = 1323.2188
$ws.Cells.Item(136, 10).Value = 2371.5557
$ws.Cells.Item(136, 11).Value = 3969.6564
$ws.Cells.Item(136, 12).Value = 7114.6671
$ws.Cells.Item(136, 13).Value = -1419.6564
$ws.Cells.Item(136, 14).Value = -12214.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 31917
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 31917
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 31917
$ws.Cells.Item(119, 14).Value = -41593

$ws.Cells.Item(122, 8).Value = 1965.174
$ws.Cells.Item(122, 9).Value = 1748.8948
$ws.Cells.Item(122, 10).Value = 2992.5
$ws.Cells.Item(122, 11).Value = 5246.6844
$ws.Cells.Item(122, 12).Value = 8977.5
$ws.Cells.Item(122, 13).Value = -2796.6844
$ws.Cells.Item(122, 14).Value = -13900

$ws.Cells.Item(126, 8).Value = 1296.6666
$ws.Cells.Item(126, 9).Value = 1330.6666
$ws.Cells.Item(126, 10).Value = 1211.6666
$ws.Cells.Item(126, 11).Value = 3991.9998
$ws.Cells.Item(126, 12).Value = 3634.9998
$ws.Cells.Item(126, 13).Value = -1521.9998
$ws.Cells.Item(126, 14).Value = -8574.9998
